# Apply updated cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.642.87"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "1.862.88"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  +1.01%  "
$ws.Range("D5").Value = "'333.44"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "'1.013"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("D7").Value = "'0.4643"
$ws.Range("E7").Value = "  -1.74%  "
$ws.Range("D8").Value = "'0.3885"
$ws.Range("E8").Value = "  -2.00%  "
$ws.Range("D9").Value = "'46.11"
$ws.Range("E9").Value = "  -3.59%  "
$ws.Range("D10").Value = "'0.07975"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").Value = "'0.9959"
$ws.Range("E11").Value = "  -3.68%  "
$ws.Range("D12").Value = "'21.51"
$ws.Range("E12").Value = "  -3.42%  "
$ws.Range("D13").Value = "1.869.20"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").Value = "'5.972"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "'7.174"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "'1.014"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").Value = "'87.81"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "'0.06721"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "'0.00001041"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "'16.88"
$ws.Range("E20").Value = "  -2.18%  "
$ws.Range("D21").Value = "'1.012"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").Value = "27.640.16"
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").Value = "'5.447"
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").Value = "'2.320"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").Value = "'158.66"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").Value = "'19.70"
$ws.Range("E27").Value = "  -2.90%  "
$ws.Range("D28").Value = "'2.105"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").Value = "'5.351"
$ws.Range("E29").Value = "  -4.64%  "
$ws.Range("D30").Value = "'121.28"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").Value = "'0.9688"
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("D32").Value = "'0.09441"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("D33").Value = "'3.646"
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("D34").Value = "'5.285"
$ws.Range("E34").Value = "  -1.90%  "
$ws.Range("D35").Value = "'1.310"
$ws.Range("E35").Value = "  -9.87%  "
$ws.Range("D36").Value = "'0.06011"
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("D37").Value = "'0.02215"
$ws.Range("D38").Value = "'1.198"
$ws.Range("E38").Value = "  -3.02%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'8.138"
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D40").Value = "'1.012"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").Value = "'0.5894"
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("D42").Value = "'0.1875"
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("D43").Value = "'10.22"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "'0.5598"
$ws.Range("E45").Value = "  -2.68%  "
$ws.Range("D46").Value = "'12.02"
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("D47").Value = "'1.912"
$ws.Range("E47").Value = "  -1.96%  "
$ws.Range("D48").Value = "'3.288"
$ws.Range("E48").Value = "  -2.80%  "
$ws.Range("D49").Value = "'0.06758"
$ws.Range("E49").Value = "  -2.44%  "
$ws.Range("D50").Value = "'112.17"
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("D51").Value = "'0.00000000297"
$ws.Range("E51").Value = "  -3.37%  "
